$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GitHub commits count (row 9)
$ws.Range("C9").Value = 38

# Update Logout score (row 31)
$ws.Range("C31").Value = 1

# Update selection / view: select B10:E10 with active cell B10
$ws.Range("B10:E10").Select()
